$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cell, [string]$text)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = $origStyle
}

Set-TextValue $ws.Range('D2') '34.287.27'
Set-TextValue $ws.Range('E2') '  +0.72%  '
Set-TextValue $ws.Range('D3') '1.794.64'
Set-TextValue $ws.Range('E3') '  +0.80%  '
Set-TextValue $ws.Range('D4') '0.998'
Set-TextValue $ws.Range('E4') '  -0.25%  '
Set-TextValue $ws.Range('D5') '227.23'
Set-TextValue $ws.Range('E5') '  +0.83%  '
Set-TextValue $ws.Range('D6') '0.549'
Set-TextValue $ws.Range('E6') '  +0.02%  '
Set-TextValue $ws.Range('D7') '0.998'
Set-TextValue $ws.Range('E7') '  -0.26%  '
Set-TextValue $ws.Range('D8') '32.39'
Set-TextValue $ws.Range('E8') '  +0.10%  '
Set-TextValue $ws.Range('D9') '0.295'
Set-TextValue $ws.Range('E9') '  +3.60%  '
Set-TextValue $ws.Range('D10') '0.0696'
Set-TextValue $ws.Range('E10') '  -1.59%  '
Set-TextValue $ws.Range('E11') '  +0.58%  '
Set-TextValue $ws.Range('D12') '2.051.77'
Set-TextValue $ws.Range('E12') '  +0.68%  '
Set-TextValue $ws.Range('D13') '11.52'
Set-TextValue $ws.Range('E13') '  +5.47%  '
Set-TextValue $ws.Range('D14') '1.781.42'
Set-TextValue $ws.Range('E14') '  -0.16%  '
Set-TextValue $ws.Range('B15') 'WrappedBTC'
Set-TextValue $ws.Range('C15') 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
Set-TextValue $ws.Range('D15') '34.207.40'
Set-TextValue $ws.Range('E15') '  +0.53%  '
Set-TextValue $ws.Range('B16') 'Polygon'
Set-TextValue $ws.Range('C16') 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
Set-TextValue $ws.Range('D16') '0.624'
Set-TextValue $ws.Range('E16') '  +1.24%  '
Set-TextValue $ws.Range('D17') '4.21'
Set-TextValue $ws.Range('E17') '  +2.10%  '
Set-TextValue $ws.Range('D18') '68.15'
Set-TextValue $ws.Range('E18') '  +0.82%  '
Set-TextValue $ws.Range('D19') '246.09'
Set-TextValue $ws.Range('E19') '  +1.12%  '
Set-TextValue $ws.Range('D20') '0.0₃0785'
Set-TextValue $ws.Range('E20') '  +0.56%  '
Set-TextValue $ws.Range('D21') '10.96'
Set-TextValue $ws.Range('E21') '  +2.65%  '
Set-TextValue $ws.Range('E22') '  -0.09%  '
Set-TextValue $ws.Range('E23') '  +1.35%  '
Set-TextValue $ws.Range('D24') '2.06'
Set-TextValue $ws.Range('E24') '  -1.01%  '
Set-TextValue $ws.Range('D25') '161.90'
Set-TextValue $ws.Range('E25') '  +1.25%  '
Set-TextValue $ws.Range('D26') '7.20'
Set-TextValue $ws.Range('E26') '  +2.67%  '
Set-TextValue $ws.Range('D27') '16.36'
Set-TextValue $ws.Range('E27') '  +0.78%  '
Set-TextValue $ws.Range('E28') '  +1.95%  '
Set-TextValue $ws.Range('E29') '  -0.10%  '
Set-TextValue $ws.Range('E30') '  +1.56%  '
Set-TextValue $ws.Range('D31') '0.0522'
Set-TextValue $ws.Range('E31') '  +2.26%  '
Set-TextValue $ws.Range('D32') '3.69'
Set-TextValue $ws.Range('E32') '  +1.82%  '
Set-TextValue $ws.Range('D33') '3.64'
Set-TextValue $ws.Range('E33') '  +4.01%  '
Set-TextValue $ws.Range('D34') '1.84'
Set-TextValue $ws.Range('E34') '  +2.05%  '
Set-TextValue $ws.Range('D35') '1.447.18'
Set-TextValue $ws.Range('E35') '  +4.18%  '
Set-TextValue $ws.Range('D36') '0.651'
Set-TextValue $ws.Range('E36') '  +0.97%  '
Set-TextValue $ws.Range('E37') '  +3.30%  '
Set-TextValue $ws.Range('D38') '2.40'
Set-TextValue $ws.Range('E38') '  +9.95%  '
Set-TextValue $ws.Range('D39') '1.04'
Set-TextValue $ws.Range('E39') '  -0.74%  '
Set-TextValue $ws.Range('D40') '81.14'
Set-TextValue $ws.Range('E40') '  +4.77%  '
Set-TextValue $ws.Range('E41') '  +2.42%  '
Set-TextValue $ws.Range('D42') '2.35'
Set-TextValue $ws.Range('E42') '  -0.01%  '
Set-TextValue $ws.Range('D43') '2.69'
Set-TextValue $ws.Range('E43') '  +0.28%  '
Set-TextValue $ws.Range('D44') '13.38'
Set-TextValue $ws.Range('E44') '  +7.81%  '
Set-TextValue $ws.Range('B45') 'BabyDogeCoin'
Set-TextValue $ws.Range('C45') 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
Set-TextValue $ws.Range('D45') '0.0₆0139'
Set-TextValue $ws.Range('E45') '  -1.53%  '
Set-TextValue $ws.Range('B46') 'FraxShare'
Set-TextValue $ws.Range('C46') 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue $ws.Range('D46') '6.08'
Set-TextValue $ws.Range('E46') '  +4.46%  '
Set-TextValue $ws.Range('B47') 'Kaspa'
Set-TextValue $ws.Range('C47') 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue $ws.Range('D47') '0.0509'
Set-TextValue $ws.Range('E47') '  +2.63%  '
Set-TextValue $ws.Range('E48') '  -0.75%  '
Set-TextValue $ws.Range('D49') '108.19'
Set-TextValue $ws.Range('E49') '  +0.73%  '
Set-TextValue $ws.Range('D50') '1.952.85'
Set-TextValue $ws.Range('E50') '  +0.72%  '
Set-TextValue $ws.Range('E51') '  -0.14%  '

Write-Host "Applied 100 cell updates"
